$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.747.30'
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Value = '3.361.85'
$ws.Range("E3").Value = '  -0.66%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.12%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '569.14'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.11%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '137.49'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.20%  '

$ws.Range("E7").Value = '  +0.00%  '

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '7.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.71%  '

$ws.Range("E10").Value = '  -1.99%  '

$ws.Range("E11").Value = '  -3.95%  '

$ws.Range("D12").Value = '3.934.65'
$ws.Range("E12").Value = '  -0.69%  '

$ws.Range("E13").Value = '  +0.60%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.89'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.81%  '

$ws.Range("D15").Value = '3.368.78'
$ws.Range("E15").Value = '  -0.46%  '

$ws.Range("E16").Value = '  -1.95%  '

$ws.Range("D17").Value = '60.872.91'
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.10'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.95%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '8.88'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -0.96%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '383.36'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.08%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.30'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +2.22%  '

$ws.Range("E23").Value = '  -2.17%  '

$ws.Range("E24").Value = '  +0.04%  '

$ws.Range("B25").Value = 'PEPE'
$ws.Range("C25").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.0000109'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -6.14%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.190'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.79%  '

$ws.Range("E27").Value = '  +0.14%  '

$ws.Range("E28").Value = '  -3.98%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.83'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -1.71%  '

$ws.Range("E30").Value = '  -1.71%  '

$ws.Range("E31").Value = '  -0.02%  '

$ws.Range("E32").Value = '  -7.16%  '

$ws.Range("E33").Value = '  -2.99%  '

$ws.Range("E34").Value = '  +0.43%  '

$ws.Range("E35").Value = '  -1.80%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.90'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -1.97%  '

$ws.Range("D37").Value = '3.395.22'
$ws.Range("E37").Value = '  -0.55%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.43'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -3.07%  '

$ws.Range("E39").Value = '  -2.72%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '25.28'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -9.03%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.768'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.52%  '

$ws.Range("E42").Value = '  -1.92%  '

$ws.Range("E43").Value = '  -3.24%  '

$ws.Range("E44").Value = '  -2.13%  '

$ws.Range("D45").Value = '2.442.07'
$ws.Range("E45").Value = '  -2.88%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.999'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.05%  '

$ws.Range("E47").Value = '  -3.25%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.06'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.31%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0257'
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  -6.13%  '

$ws.Range("E51").Value = '  -3.04%  '
